$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad) for rows 2-11 changes from serial date 45175 to 45177
# (i.e. 2023-09-06 -> 2023-09-08). Update each cell's value directly so the
# existing number format (yyyy-mm-dd style) and cell style are preserved.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45177
}
